# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Leve-profit tables across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1006182.9  # H17: 936805.8 -> 1006182.9
$ws.Cells.Item(17, 10).Value = 1006182.9  # J17: 936805.8 -> 1006182.9
$ws.Cells.Item(17, 12).Value = 3018548.7  # L17: 2810417.4 -> 3018548.7
$ws.Cells.Item(17, 14).Value = -3018884.7  # N17: -2810753.4 -> -3018884.7
$ws.Cells.Item(32, 8).Value = 1374.875  # H32: 1327.6666 -> 1374.875
$ws.Cells.Item(32, 10).Value = 1285.7142  # J32: 1243.75 -> 1285.7142
$ws.Cells.Item(32, 12).Value = 1285.7142  # L32: 1243.75 -> 1285.7142
$ws.Cells.Item(32, 14).Value = -1937.7142  # N32: -1895.75 -> -1937.7142
$ws.Cells.Item(113, 8).Value = 160945.67  # H113: 160856.78 -> 160945.67
$ws.Cells.Item(113, 9).Value = 358907.75  # I113: 287526.2 -> 358907.75
$ws.Cells.Item(113, 10).Value = 2576  # J113: 2520 -> 2576
$ws.Cells.Item(113, 11).Value = 358907.75  # K113: 287526.2 -> 358907.75
$ws.Cells.Item(113, 12).Value = 2576  # L113: 2520 -> 2576
$ws.Cells.Item(113, 13).Value = -355653.75  # M113: -284272.2 -> -355653.75
$ws.Cells.Item(113, 14).Value = -9084  # N113: -9028 -> -9084
$ws.Cells.Item(125, 8).Value = 1993.7142  # H125: 1922.2858 -> 1993.7142
$ws.Cells.Item(125, 10).Value = 858.6667  # J125: 692 -> 858.6667
$ws.Cells.Item(125, 12).Value = 7728.0003  # L125: 6228 -> 7728.0003
$ws.Cells.Item(125, 14).Value = -12648.0003  # N125: -11148 -> -12648.0003
$ws.Cells.Item(129, 8).Value = 1635.5471  # H129: 2016.9445 -> 1635.5471
$ws.Cells.Item(129, 9).Value = 691.3333  # I129: 0 -> 691.3333
$ws.Cells.Item(129, 10).Value = 1692.2  # J129: 2016.9445 -> 1692.2
$ws.Cells.Item(129, 11).Value = 2073.9999  # K129: 0 -> 2073.9999
$ws.Cells.Item(129, 12).Value = 5076.6  # L129: 6050.833500000001 -> 5076.6
$ws.Cells.Item(129, 14).Value = -15076.6  # N129: -16050.8335 -> -15076.6
$ws.Cells.Item(129, 13).Value = 2926.0001  # M129: (new) -> 2926.0001
$ws.Cells.Item(137, 8).Value = 8336727  # H137: 8003259 -> 8336727
$ws.Cells.Item(137, 9).Value = 2545.818  # I137: 2267.7693 -> 2545.818
$ws.Cells.Item(137, 10).Value = 15388727  # J137: 16671000 -> 15388727
$ws.Cells.Item(137, 11).Value = 7637.454000000001  # K137: 6803.3079 -> 7637.454000000001
$ws.Cells.Item(137, 12).Value = 46166181  # L137: 50013000 -> 46166181
$ws.Cells.Item(137, 13).Value = -5087.454000000001  # M137: -4253.3079 -> -5087.454000000001
$ws.Cells.Item(137, 14).Value = -46171281  # N137: -50018100 -> -46171281
$ws.Cells.Item(138, 8).Value = 8931071  # H138: 8624305 -> 8931071
$ws.Cells.Item(138, 9).Value = 2032.3529  # I138: 1880.1666 -> 2032.3529
$ws.Cells.Item(138, 10).Value = 22730496  # J138: 22733728 -> 22730496
$ws.Cells.Item(138, 11).Value = 6097.0587  # K138: 5640.4998 -> 6097.0587
$ws.Cells.Item(138, 12).Value = 68191488  # L138: 68201184 -> 68191488
$ws.Cells.Item(138, 13).Value = -957.0587000000005  # M138: -500.4997999999996 -> -957.0587000000005
$ws.Cells.Item(138, 14).Value = -68201768  # N138: -68211464 -> -68201768

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6917  # H32: 6681.7744 -> 6917
$ws.Cells.Item(32, 9).Value = 7097.271  # I32: 7023.9805 -> 7097.271
$ws.Cells.Item(32, 10).Value = 5955.5557  # J32: 5095.1816 -> 5955.5557
$ws.Cells.Item(32, 11).Value = 7097.271  # K32: 7023.9805 -> 7097.271
$ws.Cells.Item(32, 12).Value = 5955.5557  # L32: 5095.1816 -> 5955.5557
$ws.Cells.Item(32, 13).Value = -6810.271  # M32: -6736.9805 -> -6810.271
$ws.Cells.Item(32, 14).Value = -6529.5557  # N32: -5669.1816 -> -6529.5557
$ws.Cells.Item(61, 8).Value = 12822885  # H61: 12502413 -> 12822885
$ws.Cells.Item(61, 9).Value = 13891075  # I61: 13891156 -> 13891075
$ws.Cells.Item(61, 10).Value = 4599.3335  # J61: 3724.5 -> 4599.3335
$ws.Cells.Item(61, 11).Value = 13891075  # K61: 13891156 -> 13891075
$ws.Cells.Item(61, 12).Value = 4599.3335  # L61: 3724.5 -> 4599.3335
$ws.Cells.Item(61, 13).Value = -13890863  # M61: -13890944 -> -13890863
$ws.Cells.Item(61, 14).Value = -5023.3335  # N61: -4148.5 -> -5023.3335
$ws.Cells.Item(74, 8).Value = 9805733  # H74: 13516052 -> 9805733
$ws.Cells.Item(74, 9).Value = 17242418  # I74: 18519866 -> 17242418
$ws.Cells.Item(74, 10).Value = 2828.5908  # J74: 5753.2 -> 2828.5908
$ws.Cells.Item(74, 11).Value = 17242418  # K74: 18519866 -> 17242418
$ws.Cells.Item(74, 12).Value = 2828.5908  # L74: 5753.2 -> 2828.5908
$ws.Cells.Item(74, 13).Value = -17241544  # M74: -18518992 -> -17241544
$ws.Cells.Item(74, 14).Value = -4576.5908  # N74: -7501.2 -> -4576.5908
$ws.Cells.Item(77, 8).Value = 9805733  # H77: 13516052 -> 9805733
$ws.Cells.Item(77, 9).Value = 17242418  # I77: 18519866 -> 17242418
$ws.Cells.Item(77, 10).Value = 2828.5908  # J77: 5753.2 -> 2828.5908
$ws.Cells.Item(77, 11).Value = 86212090  # K77: 92599330 -> 86212090
$ws.Cells.Item(77, 12).Value = 14142.954  # L77: 28766 -> 14142.954
$ws.Cells.Item(77, 13).Value = -86207722  # M77: -92594962 -> -86207722
$ws.Cells.Item(77, 14).Value = -22878.954  # N77: -37502 -> -22878.954
$ws.Cells.Item(106, 8).Value = 47097.5  # H106: 47140 -> 47097.5
$ws.Cells.Item(106, 10).Value = 47097.5  # J106: 47140 -> 47097.5
$ws.Cells.Item(106, 12).Value = 47097.5  # L106: 47140 -> 47097.5
$ws.Cells.Item(106, 14).Value = -49621.5  # N106: -49664 -> -49621.5
$ws.Cells.Item(122, 8).Value = 9040  # H122: 10246.154 -> 9040
$ws.Cells.Item(122, 9).Value = 11370.182  # I122: 13630.223 -> 11370.182
$ws.Cells.Item(122, 11).Value = 34110.546  # K122: 40890.669 -> 34110.546
$ws.Cells.Item(122, 13).Value = -31660.546  # M122: -38440.669 -> -31660.546
$ws.Cells.Item(132, 8).Value = 11366943  # H132: 6946616 -> 11366943
$ws.Cells.Item(132, 9).Value = 22730438  # I132: 11906606 -> 22730438
$ws.Cells.Item(132, 10).Value = 3447.6365  # J132: 2629 -> 3447.6365
$ws.Cells.Item(132, 11).Value = 68191314  # K132: 35719818 -> 68191314
$ws.Cells.Item(132, 12).Value = 10342.9095  # L132: 7887 -> 10342.9095
$ws.Cells.Item(132, 13).Value = -68188784  # M132: -35717288 -> -68188784
$ws.Cells.Item(132, 14).Value = -15402.9095  # N132: -12947 -> -15402.9095
$ws.Cells.Item(136, 8).Value = 12822885  # H136: 12502413 -> 12822885
$ws.Cells.Item(136, 9).Value = 13891075  # I136: 13891156 -> 13891075
$ws.Cells.Item(136, 10).Value = 4599.3335  # J136: 3724.5 -> 4599.3335
$ws.Cells.Item(136, 11).Value = 41673225  # K136: 41673468 -> 41673225
$ws.Cells.Item(136, 12).Value = 13798.0005  # L136: 11173.5 -> 13798.0005
$ws.Cells.Item(136, 13).Value = -41670675  # M136: -41670918 -> -41670675
$ws.Cells.Item(136, 14).Value = -18898.0005  # N136: -16273.5 -> -18898.0005

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4221.511  # H134: 4009.75 -> 4221.511
$ws.Cells.Item(134, 9).Value = 3175.04  # I134: 2764.6897 -> 3175.04
$ws.Cells.Item(134, 10).Value = 5529.6  # J134: 5910.1055 -> 5529.6
$ws.Cells.Item(134, 11).Value = 9525.119999999999  # K134: 8294.069100000001 -> 9525.119999999999
$ws.Cells.Item(134, 12).Value = 16588.8  # L134: 17730.3165 -> 16588.8
$ws.Cells.Item(134, 13).Value = -6990.119999999999  # M134: -5759.069100000001 -> -6990.119999999999
$ws.Cells.Item(134, 14).Value = -21658.8  # N134: -22800.3165 -> -21658.8

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(21, 8).Value = 6100  # H21: 0 -> 6100
$ws.Cells.Item(21, 10).Value = 6100  # J21: 0 -> 6100
$ws.Cells.Item(21, 12).Value = 6100  # L21: 0 -> 6100
$ws.Cells.Item(21, 14).Value = -6570  # N21: (new) -> -6570
$ws.Cells.Item(56, 8).Value = 0  # H56: 10 -> 0
$ws.Cells.Item(56, 9).Value = 0  # I56: 10 -> 0
$ws.Cells.Item(56, 11).Value = 0  # K56: 10 -> 0
$ws.Cells.Item(56, 13).ClearContents()  # M56: removed
$ws.Cells.Item(58, 8).Value = 2417.1724  # H58: 2493.8928 -> 2417.1724
$ws.Cells.Item(58, 9).Value = 1816.3846  # I58: 1826.4615 -> 1816.3846
$ws.Cells.Item(58, 10).Value = 2905.3125  # J58: 3072.3333 -> 2905.3125
$ws.Cells.Item(58, 11).Value = 1816.3846  # K58: 1826.4615 -> 1816.3846
$ws.Cells.Item(58, 12).Value = 2905.3125  # L58: 3072.3333 -> 2905.3125
$ws.Cells.Item(58, 13).Value = -1613.3846  # M58: -1623.4615 -> -1613.3846
$ws.Cells.Item(58, 14).Value = -3311.3125  # N58: -3478.3333 -> -3311.3125
$ws.Cells.Item(132, 8).Value = 3839.7334  # H132: 4954.1816 -> 3839.7334
$ws.Cells.Item(132, 9).Value = 3993.7144  # I132: 5271.2 -> 3993.7144
$ws.Cells.Item(132, 10).Value = 3705  # J132: 4690 -> 3705
$ws.Cells.Item(132, 11).Value = 11981.1432  # K132: 15813.6 -> 11981.1432
$ws.Cells.Item(132, 12).Value = 11115  # L132: 14070 -> 11115
$ws.Cells.Item(132, 13).Value = -9451.143199999999  # M132: -13283.6 -> -9451.143199999999
$ws.Cells.Item(132, 14).Value = -16175  # N132: -19130 -> -16175
$ws.Cells.Item(134, 8).Value = 884622.9399999999  # H134: 955329.2 -> 884622.9399999999
$ws.Cells.Item(134, 9).Value = 3556.375  # I134: 3582 -> 3556.375
$ws.Cells.Item(134, 10).Value = 2166174.2  # J134: 2647324.2 -> 2166174.2
$ws.Cells.Item(134, 11).Value = 10669.125  # K134: 10746 -> 10669.125
$ws.Cells.Item(134, 12).Value = 6498522.600000001  # L134: 7941972.600000001 -> 6498522.600000001
$ws.Cells.Item(134, 13).Value = -8134.125  # M134: -8211 -> -8134.125
$ws.Cells.Item(134, 14).Value = -6503592.600000001  # N134: -7947042.600000001 -> -6503592.600000001
$ws.Cells.Item(136, 8).Value = 2417.1724  # H136: 2493.8928 -> 2417.1724
$ws.Cells.Item(136, 9).Value = 1816.3846  # I136: 1826.4615 -> 1816.3846
$ws.Cells.Item(136, 10).Value = 2905.3125  # J136: 3072.3333 -> 2905.3125
$ws.Cells.Item(136, 11).Value = 5449.1538  # K136: 5479.3845 -> 5449.1538
$ws.Cells.Item(136, 12).Value = 8715.9375  # L136: 9216.999899999999 -> 8715.9375
$ws.Cells.Item(136, 13).Value = -2899.1538  # M136: -2929.3845 -> -2899.1538
$ws.Cells.Item(136, 14).Value = -13815.9375  # N136: -14316.9999 -> -13815.9375

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 98.42856999999999  # H23: 65.69231000000001 -> 98.42856999999999
$ws.Cells.Item(23, 9).Value = 57.8  # I23: 75.5 -> 57.8
$ws.Cells.Item(23, 10).Value = 200  # J23: 63.909092 -> 200
$ws.Cells.Item(23, 11).Value = 173.4  # K23: 226.5 -> 173.4
$ws.Cells.Item(23, 12).Value = 600  # L23: 191.727276 -> 600
$ws.Cells.Item(23, 13).Value = 61.60000000000002  # M23: 8.5 -> 61.60000000000002
$ws.Cells.Item(23, 14).Value = -1070  # N23: -661.7272760000001 -> -1070
$ws.Cells.Item(56, 8).Value = 4429.1665  # H56: 4608 -> 4429.1665
$ws.Cells.Item(56, 9).Value = 4429.1665  # I56: 4608 -> 4429.1665
$ws.Cells.Item(56, 11).Value = 4429.1665  # K56: 4608 -> 4429.1665
$ws.Cells.Item(56, 13).Value = -3899.1665  # M56: -4078 -> -3899.1665
$ws.Cells.Item(68, 8).Value = 1304.3695  # H68: 1325.25 -> 1304.3695
$ws.Cells.Item(68, 9).Value = 1028.8334  # I68: 1048.1765 -> 1028.8334
$ws.Cells.Item(68, 10).Value = 1481.5  # J68: 1499.7037 -> 1481.5
$ws.Cells.Item(68, 11).Value = 3086.5002  # K68: 3144.5295 -> 3086.5002
$ws.Cells.Item(68, 12).Value = 4444.5  # L68: 4499.1111 -> 4444.5
$ws.Cells.Item(68, 13).Value = -2275.5002  # M68: -2333.5295 -> -2275.5002
$ws.Cells.Item(68, 14).Value = -6066.5  # N68: -6121.1111 -> -6066.5
$ws.Cells.Item(71, 8).Value = 1304.3695  # H71: 1325.25 -> 1304.3695
$ws.Cells.Item(71, 9).Value = 1028.8334  # I71: 1048.1765 -> 1028.8334
$ws.Cells.Item(71, 10).Value = 1481.5  # J71: 1499.7037 -> 1481.5
$ws.Cells.Item(71, 11).Value = 9259.500599999999  # K71: 9433.5885 -> 9259.500599999999
$ws.Cells.Item(71, 12).Value = 13333.5  # L71: 13497.3333 -> 13333.5
$ws.Cells.Item(71, 13).Value = -5203.500599999999  # M71: -5377.5885 -> -5203.500599999999
$ws.Cells.Item(71, 14).Value = -21445.5  # N71: -21609.3333 -> -21445.5
$ws.Cells.Item(86, 8).Value = 1000  # H86: 1347.2858 -> 1000
$ws.Cells.Item(86, 9).Value = 1000  # I86: 816 -> 1000
$ws.Cells.Item(86, 10).Value = 1000  # J86: 1559.8 -> 1000
$ws.Cells.Item(86, 11).Value = 3000  # K86: 2448 -> 3000
$ws.Cells.Item(86, 12).Value = 3000  # L86: 4679.4 -> 3000
$ws.Cells.Item(86, 13).Value = -1814  # M86: -1262 -> -1814
$ws.Cells.Item(86, 14).Value = -5372  # N86: -7051.4 -> -5372
$ws.Cells.Item(89, 8).Value = 1000  # H89: 1347.2858 -> 1000
$ws.Cells.Item(89, 9).Value = 1000  # I89: 816 -> 1000
$ws.Cells.Item(89, 10).Value = 1000  # J89: 1559.8 -> 1000
$ws.Cells.Item(89, 11).Value = 9000  # K89: 7344 -> 9000
$ws.Cells.Item(89, 12).Value = 9000  # L89: 14038.2 -> 9000
$ws.Cells.Item(89, 13).Value = -3072  # M89: -1416 -> -3072
$ws.Cells.Item(89, 14).Value = -20856  # N89: -25894.2 -> -20856
$ws.Cells.Item(107, 8).Value = 973.913  # H107: 975.0625 -> 973.913
$ws.Cells.Item(107, 10).Value = 1582.25  # J107: 1537.5769 -> 1582.25
$ws.Cells.Item(107, 12).Value = 4746.75  # L107: 4612.7307 -> 4746.75
$ws.Cells.Item(107, 14).Value = -8586.75  # N107: -8452.7307 -> -8586.75
$ws.Cells.Item(131, 8).Value = 740.46  # H131: 746.66 -> 740.46
$ws.Cells.Item(131, 9).Value = 390.8421  # I131: 428.73685 -> 390.8421
$ws.Cells.Item(131, 10).Value = 822.4691  # J131: 821.23456 -> 822.4691
$ws.Cells.Item(131, 11).Value = 1172.5263  # K131: 1286.21055 -> 1172.5263
$ws.Cells.Item(131, 12).Value = 2467.4073  # L131: 2463.70368 -> 2467.4073
$ws.Cells.Item(131, 13).Value = 3867.4737  # M131: 3753.78945 -> 3867.4737
$ws.Cells.Item(131, 14).Value = -12547.4073  # N131: -12543.70368 -> -12547.4073

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 4169391.5  # H122: 5558914 -> 4169391.5
$ws.Cells.Item(122, 9).Value = 5130674  # I122: 9526881 -> 5130674
$ws.Cells.Item(122, 10).Value = 3833.3333  # J122: 3760 -> 3833.3333
$ws.Cells.Item(122, 11).Value = 15392022  # K122: 28580643 -> 15392022
$ws.Cells.Item(122, 12).Value = 11499.9999  # L122: 11280 -> 11499.9999
$ws.Cells.Item(122, 13).Value = -15389572  # M122: -28578193 -> -15389572
$ws.Cells.Item(122, 14).Value = -16399.9999  # N122: -16180 -> -16399.9999
$ws.Cells.Item(126, 8).Value = 4247.3687  # H126: 4729.353 -> 4247.3687
$ws.Cells.Item(126, 9).Value = 2557.1428  # I126: 2649.8333 -> 2557.1428
$ws.Cells.Item(126, 10).Value = 5233.3335  # J126: 5863.636 -> 5233.3335
$ws.Cells.Item(126, 11).Value = 7671.428400000001  # K126: 7949.499899999999 -> 7671.428400000001
$ws.Cells.Item(126, 12).Value = 15700.0005  # L126: 17590.908 -> 15700.0005
$ws.Cells.Item(126, 13).Value = -5201.428400000001  # M126: -5479.499899999999 -> -5201.428400000001
$ws.Cells.Item(126, 14).Value = -20640.0005  # N126: -22530.908 -> -20640.0005
$ws.Cells.Item(132, 8).Value = 5892.5186  # H132: 5448.3105 -> 5892.5186
$ws.Cells.Item(132, 9).Value = 5138.1333  # I132: 5151.467 -> 5138.1333
$ws.Cells.Item(132, 10).Value = 6835.5  # J132: 5766.357 -> 6835.5
$ws.Cells.Item(132, 11).Value = 15414.3999  # K132: 15454.401 -> 15414.3999
$ws.Cells.Item(132, 12).Value = 20506.5  # L132: 17299.071 -> 20506.5
$ws.Cells.Item(132, 13).Value = -12884.3999  # M132: -12924.401 -> -12884.3999
$ws.Cells.Item(132, 14).Value = -25566.5  # N132: -22359.071 -> -25566.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 12833.333  # H40: 16400 -> 12833.333
$ws.Cells.Item(40, 9).Value = 18333.334  # I40: 26500 -> 18333.334
$ws.Cells.Item(40, 10).Value = 7333.3335  # J40: 9666.666999999999 -> 7333.3335
$ws.Cells.Item(40, 11).Value = 18333.334  # K40: 26500 -> 18333.334
$ws.Cells.Item(40, 12).Value = 7333.3335  # L40: 9666.666999999999 -> 7333.3335
$ws.Cells.Item(40, 13).Value = -18197.334  # M40: -26364 -> -18197.334
$ws.Cells.Item(40, 14).Value = -7605.3335  # N40: -9938.666999999999 -> -7605.3335
$ws.Cells.Item(56, 8).Value = 39800  # H56: 25600 -> 39800
$ws.Cells.Item(56, 9).Value = 0  # I56: 19800 -> 0
$ws.Cells.Item(56, 10).Value = 39800  # J56: 31400 -> 39800
$ws.Cells.Item(56, 11).Value = 0  # K56: 19800 -> 0
$ws.Cells.Item(56, 14).Value = -41182  # N56: -32782 -> -41182
$ws.Cells.Item(56, 13).ClearContents()  # M56: removed
$ws.Cells.Item(68, 8).Value = 2180.2  # H68: 2080 -> 2180.2
$ws.Cells.Item(68, 10).Value = 2350.25  # J68: 2196 -> 2350.25
$ws.Cells.Item(68, 12).Value = 2350.25  # L68: 2196 -> 2350.25
$ws.Cells.Item(68, 14).Value = -3848.25  # N68: -3694 -> -3848.25
$ws.Cells.Item(71, 8).Value = 2180.2  # H71: 2080 -> 2180.2
$ws.Cells.Item(71, 10).Value = 2350.25  # J71: 2196 -> 2350.25
$ws.Cells.Item(71, 12).Value = 11751.25  # L71: 10980 -> 11751.25
$ws.Cells.Item(71, 14).Value = -19239.25  # N71: -18468 -> -19239.25
$ws.Cells.Item(106, 8).Value = 14048.889  # H106: 20300 -> 14048.889
$ws.Cells.Item(106, 10).Value = 14048.889  # J106: 20300 -> 14048.889
$ws.Cells.Item(106, 12).Value = 14048.889  # L106: 20300 -> 14048.889
$ws.Cells.Item(106, 14).Value = -16572.889  # N106: -22824 -> -16572.889
$ws.Cells.Item(136, 8).Value = 53577930  # H136: 11115438 -> 53577930
$ws.Cells.Item(136, 9).Value = 88236870  # I136: 21740450 -> 88236870
$ws.Cells.Item(136, 10).Value = 14102.272  # J136: 7471.5454 -> 14102.272
$ws.Cells.Item(136, 11).Value = 264710610  # K136: 65221350 -> 264710610
$ws.Cells.Item(136, 12).Value = 42306.81600000001  # L136: 22414.6362 -> 42306.81600000001
$ws.Cells.Item(136, 13).Value = -264708060  # M136: -65218800 -> -264708060
$ws.Cells.Item(136, 14).Value = -47406.81600000001  # N136: -27514.6362 -> -47406.81600000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2431.8333  # H81: 2552.3845 -> 2431.8333
$ws.Cells.Item(81, 9).Value = 1882.8572  # I81: 2084.875 -> 1882.8572
$ws.Cells.Item(81, 10).Value = 3200.4  # J81: 3300.4 -> 3200.4
$ws.Cells.Item(81, 11).Value = 3765.7144  # K81: 4169.75 -> 3765.7144
$ws.Cells.Item(81, 12).Value = 6400.8  # L81: 6600.8 -> 6400.8
$ws.Cells.Item(81, 13).Value = -2704.7144  # M81: -3108.75 -> -2704.7144
$ws.Cells.Item(81, 14).Value = -8522.799999999999  # N81: -8722.799999999999 -> -8522.799999999999
$ws.Cells.Item(84, 8).Value = 2431.8333  # H84: 2552.3845 -> 2431.8333
$ws.Cells.Item(84, 9).Value = 1882.8572  # I84: 2084.875 -> 1882.8572
$ws.Cells.Item(84, 10).Value = 3200.4  # J84: 3300.4 -> 3200.4
$ws.Cells.Item(84, 11).Value = 18828.572  # K84: 20848.75 -> 18828.572
$ws.Cells.Item(84, 12).Value = 32004  # L84: 33004 -> 32004
$ws.Cells.Item(84, 13).Value = -13524.572  # M84: -15544.75 -> -13524.572
$ws.Cells.Item(84, 14).Value = -42612  # N84: -43612 -> -42612
$ws.Cells.Item(107, 8).Value = 1617.8334  # H107: 1746.5454 -> 1617.8334
$ws.Cells.Item(107, 9).Value = 2035.5  # I107: 2311.7144 -> 2035.5
$ws.Cells.Item(107, 10).Value = 782.5  # J107: 757.5 -> 782.5
$ws.Cells.Item(107, 11).Value = 6106.5  # K107: 6935.1432 -> 6106.5
$ws.Cells.Item(107, 12).Value = 2347.5  # L107: 2272.5 -> 2347.5
$ws.Cells.Item(107, 13).Value = -4186.5  # M107: -5015.1432 -> -4186.5
$ws.Cells.Item(107, 14).Value = -6187.5  # N107: -6112.5 -> -6187.5
$ws.Cells.Item(122, 8).Value = 2538.4443  # H122: 2712.5334 -> 2538.4443
$ws.Cells.Item(122, 9).Value = 2292.8  # I122: 2491.3845 -> 2292.8
$ws.Cells.Item(122, 10).Value = 3766.6667  # J122: 4150 -> 3766.6667
$ws.Cells.Item(122, 11).Value = 6878.400000000001  # K122: 7474.1535 -> 6878.400000000001
$ws.Cells.Item(122, 12).Value = 11300.0001  # L122: 12450 -> 11300.0001
$ws.Cells.Item(122, 13).Value = -4428.400000000001  # M122: -5024.1535 -> -4428.400000000001
$ws.Cells.Item(122, 14).Value = -16200.0001  # N122: -17350 -> -16200.0001
$ws.Cells.Item(132, 8).Value = 6922.048  # H132: 7587.5264 -> 6922.048
$ws.Cells.Item(132, 9).Value = 9301.071  # I132: 9978.076999999999 -> 9301.071
$ws.Cells.Item(132, 10).Value = 2164  # J132: 2408 -> 2164
$ws.Cells.Item(132, 11).Value = 27903.213  # K132: 29934.231 -> 27903.213
$ws.Cells.Item(132, 12).Value = 6492  # L132: 7224 -> 6492
$ws.Cells.Item(132, 13).Value = -25373.213  # M132: -27404.231 -> -25373.213
$ws.Cells.Item(132, 14).Value = -11552  # N132: -12284 -> -11552
$ws.Cells.Item(136, 8).Value = 10001790  # H136: 3334265.5 -> 10001790
$ws.Cells.Item(136, 9).Value = 11112988  # I136: 3572320 -> 11112988
$ws.Cells.Item(136, 10).Value = 1000  # J136: 1502.5 -> 1000
$ws.Cells.Item(136, 11).Value = 33338964  # K136: 10716960 -> 33338964
$ws.Cells.Item(136, 12).Value = 3000  # L136: 4507.5 -> 3000
$ws.Cells.Item(136, 13).Value = -33336414  # M136: -10714410 -> -33336414
$ws.Cells.Item(136, 14).Value = -8100  # N136: -9607.5 -> -8100

Write-Output "Applied Ultima_Profits updates to all sheets"
